# BM Table.xlsx edit: populate GSSA sheet with simulation results, wire the
# Table sheet's "GSSA" column to them, and add the two new RMS-Euler-error
# summary figures (D9, D11, D12) on the Table sheet.

$wb = $excel.ActiveWorkbook

$tblws  = $wb.Worksheets.Item("Table")
$gssaws = $wb.Worksheets.Item("GSSA")

# --- GSSA sheet: fill in the TZFE/OPFE mean & RMSE block (B3:E10) ---------
$gssaws.Range("B3:E10").NumberFormat = "0.00%"

$gssaData = @(
    @(0.0290809,    0.40088,     0.0682931,   0.507325),
    @(0.00961876,   0.132277,    0.00725493,  0.00909414),
    @(0.0155769,    0.214213,    0.0117489,   0.0147274),
    @(0.000112862,  0.0105704,   0.0039484,   0.00494937),
    @(0.192375,     2.64553,     0.145099,    0.181883),
    @(0.0154935,    0.213066,    0.011686,    0.0146485),
    @(0.0292208,    0.401843,    0.00826206,  0.0104971),
    @(0.00454358,   0.0650114,   0.00128636,  0.00163617)
)

for ($i = 0; $i -lt $gssaData.Length; $i++) {
    $row = 3 + $i
    $vals = $gssaData[$i]
    $gssaws.Cells.Item($row, 2).Value = $vals[0]
    $gssaws.Cells.Item($row, 3).Value = $vals[1]
    $gssaws.Cells.Item($row, 4).Value = $vals[2]
    $gssaws.Cells.Item($row, 5).Value = $vals[3]
}

# --- Table sheet: point the (previously empty) GSSA column at GSSA sheet -
$tblws.Range("D3").Formula = "=GSSA!B3"
$tblws.Range("D4").Formula = "=GSSA!C3"
$tblws.Range("D6").Formula = "=GSSA!D3"
$tblws.Range("D7").Formula = "=GSSA!E3"

# --- Table sheet: new RMS Euler-error summary figures ---------------------
$tblws.Range("D9").NumberFormat = "0.00E+00"
$tblws.Range("D9").Value = 0.000241043

$tblws.Range("D11").Value = 173.34756409923267
$tblws.Range("D12").Value = 122.76904415020975

$excel.Calculate()

# --- Restore cursor/selection on each touched sheet, Table active last ---
$gssaws.Activate()
$null = $gssaws.Range("E14").Select()

$null = $tblws.Activate()
$null = $tblws.Range("D12").Select()
